$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": update rows 2-34 ---
$ws1.Cells.Item(2, 1).Value = "BRVM - CONSOMMATION DE BASE     (**)"
$ws1.Cells.Item(2, 2).Value = 0
$ws1.Cells.Item(2, 3).Value = 4
$ws1.Cells.Item(2, 4).Value = 1059.31
$ws1.Cells.Item(2, 5).Value = 273.92
$ws1.Cells.Item(2, 6).Value = "🟡 Observer"
$ws1.Cells.Item(2, 7).Value = "➖ Neutre"

$ws1.Cells.Item(3, 1).Value = "BRVM-PRINCIPAL     (**)"
$ws1.Cells.Item(3, 2).Value = 0
$ws1.Cells.Item(3, 3).Value = 4
$ws1.Cells.Item(3, 4).Value = 1023.24
$ws1.Cells.Item(3, 5).Value = 260.59
$ws1.Cells.Item(3, 6).Value = "🟡 Observer"
$ws1.Cells.Item(3, 7).Value = "➖ Neutre"

$ws1.Cells.Item(4, 1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(4, 2).Value = 0
$ws1.Cells.Item(4, 3).Value = 4
$ws1.Cells.Item(4, 4).Value = 774.29
$ws1.Cells.Item(4, 5).Value = 185.83
$ws1.Cells.Item(4, 6).Value = "🟡 Observer"
$ws1.Cells.Item(4, 7).Value = "➖ Neutre"

$ws1.Cells.Item(5, 1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(5, 2).Value = 0
$ws1.Cells.Item(5, 3).Value = 4
$ws1.Cells.Item(5, 4).Value = 765.48
$ws1.Cells.Item(5, 5).Value = 188.28
$ws1.Cells.Item(5, 6).Value = "🟡 Observer"
$ws1.Cells.Item(5, 7).Value = "➖ Neutre"

$ws1.Cells.Item(6, 1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(6, 2).Value = 0
$ws1.Cells.Item(6, 3).Value = 4
$ws1.Cells.Item(6, 4).Value = 641.36
$ws1.Cells.Item(6, 5).Value = 162.17
$ws1.Cells.Item(6, 6).Value = "🟡 Observer"
$ws1.Cells.Item(6, 7).Value = "➖ Neutre"

$ws1.Cells.Item(7, 1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(7, 2).Value = 0
$ws1.Cells.Item(7, 3).Value = 4
$ws1.Cells.Item(7, 4).Value = 617.18
$ws1.Cells.Item(7, 5).Value = 155.31
$ws1.Cells.Item(7, 6).Value = "🟡 Observer"
$ws1.Cells.Item(7, 7).Value = "➖ Neutre"

$ws1.Cells.Item(8, 1).Value = "BRVM – COMPOSITE TOTAL RETURN     (**)"
$ws1.Cells.Item(8, 2).Value = 0
$ws1.Cells.Item(8, 3).Value = 4
$ws1.Cells.Item(8, 4).Value = 589.93
$ws1.Cells.Item(8, 5).Value = 149.31
$ws1.Cells.Item(8, 6).Value = "🟡 Observer"
$ws1.Cells.Item(8, 7).Value = "➖ Neutre"

$ws1.Cells.Item(9, 1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(9, 2).Value = 0
$ws1.Cells.Item(9, 3).Value = 4
$ws1.Cells.Item(9, 4).Value = 502.27
$ws1.Cells.Item(9, 5).Value = 127.31
$ws1.Cells.Item(9, 6).Value = "🟡 Observer"
$ws1.Cells.Item(9, 7).Value = "➖ Neutre"

$ws1.Cells.Item(10, 1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(10, 2).Value = 0
$ws1.Cells.Item(10, 3).Value = 4
$ws1.Cells.Item(10, 4).Value = 477.05
$ws1.Cells.Item(10, 5).Value = 119.62
$ws1.Cells.Item(10, 6).Value = "🟡 Observer"
$ws1.Cells.Item(10, 7).Value = "➖ Neutre"

$ws1.Cells.Item(11, 1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(11, 2).Value = 0
$ws1.Cells.Item(11, 3).Value = 4
$ws1.Cells.Item(11, 4).Value = 400.87
$ws1.Cells.Item(11, 5).Value = 101.15
$ws1.Cells.Item(11, 6).Value = "🟡 Observer"
$ws1.Cells.Item(11, 7).Value = "➖ Neutre"

$ws1.Cells.Item(12, 1).Value = "UNILEVER CI (UNLC)"
$ws1.Cells.Item(12, 2).Value = 3
$ws1.Cells.Item(12, 3).Value = 0
$ws1.Cells.Item(12, 4).Value = 20.44
$ws1.Cells.Item(12, 5).Value = 5.53
$ws1.Cells.Item(12, 6).Value = "🟢 Achat"
$ws1.Cells.Item(12, 7).Value = "✅ Renforcer"

$ws1.Cells.Item(13, 1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(13, 2).Value = 3
$ws1.Cells.Item(13, 3).Value = 1
$ws1.Cells.Item(13, 4).Value = 11.42
$ws1.Cells.Item(13, 5).Value = 4.27
$ws1.Cells.Item(13, 6).Value = "🟢 Achat"
$ws1.Cells.Item(13, 7).Value = "✅ Renforcer"

$ws1.Cells.Item(14, 1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(14, 2).Value = 2
$ws1.Cells.Item(14, 3).Value = 0
$ws1.Cells.Item(14, 4).Value = 11.26
$ws1.Cells.Item(14, 5).Value = 7.41
$ws1.Cells.Item(14, 6).Value = "🟡 Observer"
$ws1.Cells.Item(14, 7).Value = "➖ Neutre"

$ws1.Cells.Item(15, 1).Value = "SETAO CI (STAC)"
$ws1.Cells.Item(15, 2).Value = 2
$ws1.Cells.Item(15, 3).Value = 1
$ws1.Cells.Item(15, 4).Value = 7.73
$ws1.Cells.Item(15, 5).Value = -6.83
$ws1.Cells.Item(15, 6).Value = "🟡 Observer"
$ws1.Cells.Item(15, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(16, 1).Value = "SMB CI (SMBC)"
$ws1.Cells.Item(16, 2).Value = 1
$ws1.Cells.Item(16, 3).Value = 0
$ws1.Cells.Item(16, 4).Value = 7.49
$ws1.Cells.Item(16, 5).Value = 7.49
$ws1.Cells.Item(16, 6).Value = "🟡 Observer"
$ws1.Cells.Item(16, 7).Value = "➖ Neutre"

$ws1.Cells.Item(17, 1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(17, 2).Value = 1
$ws1.Cells.Item(17, 3).Value = 0
$ws1.Cells.Item(17, 4).Value = 7.45
$ws1.Cells.Item(17, 5).Value = 7.45
$ws1.Cells.Item(17, 6).Value = "🟡 Observer"
$ws1.Cells.Item(17, 7).Value = "➖ Neutre"

$ws1.Cells.Item(18, 1).Value = "SICOR CI (SICC)"
$ws1.Cells.Item(18, 2).Value = 1
$ws1.Cells.Item(18, 3).Value = 0
$ws1.Cells.Item(18, 4).Value = 7.44
$ws1.Cells.Item(18, 5).Value = 7.44
$ws1.Cells.Item(18, 6).Value = "🟡 Observer"
$ws1.Cells.Item(18, 7).Value = "➖ Neutre"

$ws1.Cells.Item(19, 1).Value = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(19, 2).Value = 1
$ws1.Cells.Item(19, 3).Value = 0
$ws1.Cells.Item(19, 4).Value = 2.17
$ws1.Cells.Item(19, 5).Value = 2.17
$ws1.Cells.Item(19, 6).Value = "🟡 Observer"
$ws1.Cells.Item(19, 7).Value = "➖ Neutre"

$ws1.Cells.Item(20, 1).Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Cells.Item(20, 2).Value = 1
$ws1.Cells.Item(20, 3).Value = 0
$ws1.Cells.Item(20, 4).Value = 1.65
$ws1.Cells.Item(20, 5).Value = 1.65
$ws1.Cells.Item(20, 6).Value = "🟡 Observer"
$ws1.Cells.Item(20, 7).Value = "➖ Neutre"

$ws1.Cells.Item(21, 1).Value = "EVIOSYS PACKAGING SIEM CI (SEMC)"
$ws1.Cells.Item(21, 2).Value = 2
$ws1.Cells.Item(21, 3).Value = 2
$ws1.Cells.Item(21, 4).Value = 0.01
$ws1.Cells.Item(21, 5).Value = -7.36
$ws1.Cells.Item(21, 6).Value = "🟡 Observer"
$ws1.Cells.Item(21, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(22, 1).Value = "ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)"
$ws1.Cells.Item(22, 2).Value = 2
$ws1.Cells.Item(22, 3).Value = 2
$ws1.Cells.Item(22, 4).Value = -0.36
$ws1.Cells.Item(22, 5).Value = -7.39
$ws1.Cells.Item(22, 6).Value = "🟡 Observer"
$ws1.Cells.Item(22, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(23, 1).Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws1.Cells.Item(23, 2).Value = 0
$ws1.Cells.Item(23, 3).Value = 1
$ws1.Cells.Item(23, 4).Value = -1.59
$ws1.Cells.Item(23, 5).Value = -1.59
$ws1.Cells.Item(23, 6).Value = "🟡 Observer"
$ws1.Cells.Item(23, 7).Value = "➖ Neutre"

$ws1.Cells.Item(24, 1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Cells.Item(24, 2).Value = 0
$ws1.Cells.Item(24, 3).Value = 1
$ws1.Cells.Item(24, 4).Value = -1.64
$ws1.Cells.Item(24, 5).Value = -1.64
$ws1.Cells.Item(24, 6).Value = "🟡 Observer"
$ws1.Cells.Item(24, 7).Value = "➖ Neutre"

$ws1.Cells.Item(25, 1).Value = "CORIS BANK INTERNATIONAL (CBIBF)"
$ws1.Cells.Item(25, 2).Value = 0
$ws1.Cells.Item(25, 3).Value = 1
$ws1.Cells.Item(25, 4).Value = -2
$ws1.Cells.Item(25, 5).Value = -2
$ws1.Cells.Item(25, 6).Value = "🟡 Observer"
$ws1.Cells.Item(25, 7).Value = "➖ Neutre"

$ws1.Cells.Item(26, 1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Cells.Item(26, 2).Value = 0
$ws1.Cells.Item(26, 3).Value = 1
$ws1.Cells.Item(26, 4).Value = -2.29
$ws1.Cells.Item(26, 5).Value = -2.29
$ws1.Cells.Item(26, 6).Value = "🟡 Observer"
$ws1.Cells.Item(26, 7).Value = "➖ Neutre"

$ws1.Cells.Item(27, 1).Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Cells.Item(27, 2).Value = 0
$ws1.Cells.Item(27, 3).Value = 1
$ws1.Cells.Item(27, 4).Value = -2.72
$ws1.Cells.Item(27, 5).Value = -2.72
$ws1.Cells.Item(27, 6).Value = "🟡 Observer"
$ws1.Cells.Item(27, 7).Value = "➖ Neutre"

$ws1.Cells.Item(28, 1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(28, 2).Value = 0
$ws1.Cells.Item(28, 3).Value = 1
$ws1.Cells.Item(28, 4).Value = -2.78
$ws1.Cells.Item(28, 5).Value = -2.78
$ws1.Cells.Item(28, 6).Value = "🟡 Observer"
$ws1.Cells.Item(28, 7).Value = "➖ Neutre"

$ws1.Cells.Item(29, 1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(29, 2).Value = 0
$ws1.Cells.Item(29, 3).Value = 1
$ws1.Cells.Item(29, 4).Value = -2.94
$ws1.Cells.Item(29, 5).Value = -2.94
$ws1.Cells.Item(29, 6).Value = "🟡 Observer"
$ws1.Cells.Item(29, 7).Value = "➖ Neutre"

$ws1.Cells.Item(30, 1).Value = "ONATEL BF (ONTBF)"
$ws1.Cells.Item(30, 2).Value = 0
$ws1.Cells.Item(30, 3).Value = 1
$ws1.Cells.Item(30, 4).Value = -2.95
$ws1.Cells.Item(30, 5).Value = -2.95
$ws1.Cells.Item(30, 6).Value = "🟡 Observer"
$ws1.Cells.Item(30, 7).Value = "➖ Neutre"

$ws1.Cells.Item(31, 1).Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Cells.Item(31, 2).Value = 1
$ws1.Cells.Item(31, 3).Value = 2
$ws1.Cells.Item(31, 4).Value = -3.98
$ws1.Cells.Item(31, 5).Value = 1.79
$ws1.Cells.Item(31, 6).Value = "🟡 Observer"
$ws1.Cells.Item(31, 7).Value = "👀 À surveiller"

$ws1.Cells.Item(32, 1).Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws1.Cells.Item(32, 2).Value = 0
$ws1.Cells.Item(32, 3).Value = 1
$ws1.Cells.Item(32, 4).Value = -4.67
$ws1.Cells.Item(32, 5).Value = -4.67
$ws1.Cells.Item(32, 6).Value = "🟡 Observer"
$ws1.Cells.Item(32, 7).Value = "➖ Neutre"

$ws1.Cells.Item(33, 1).Value = "FILTISAC CI (FTSC)"
$ws1.Cells.Item(33, 2).Value = 0
$ws1.Cells.Item(33, 3).Value = 1
$ws1.Cells.Item(33, 4).Value = -7.26
$ws1.Cells.Item(33, 5).Value = -7.26
$ws1.Cells.Item(33, 6).Value = "🟡 Observer"
$ws1.Cells.Item(33, 7).Value = "➖ Neutre"

$ws1.Cells.Item(34, 1).Value = "SICABLE CI (CABC)"
$ws1.Cells.Item(34, 2).Value = 0
$ws1.Cells.Item(34, 3).Value = 2
$ws1.Cells.Item(34, 4).Value = -14.93
$ws1.Cells.Item(34, 5).Value = -7.43
$ws1.Cells.Item(34, 6).Value = "🟡 Observer"
$ws1.Cells.Item(34, 7).Value = "➖ Neutre"

# --- Sheet "Top_YTD": update rows 2-11 (B column only; titles unchanged) ---
$ws2.Cells.Item(2, 2).Value = 17603.87
$ws2.Cells.Item(3, 2).Value = 15924.36
$ws2.Cells.Item(4, 2).Value = 7323.08
$ws2.Cells.Item(5, 2).Value = 7106.69
$ws2.Cells.Item(6, 2).Value = 4493.51
$ws2.Cells.Item(7, 2).Value = 4081.64
$ws2.Cells.Item(8, 2).Value = 3651.09
$ws2.Cells.Item(9, 2).Value = 2488.7
$ws2.Cells.Item(10, 2).Value = 2211.24
$ws2.Cells.Item(11, 2).Value = 1506.94
